# Move the "4875 Primrose Ln" / "Livermore, CA 94551" address paragraphs so
# that they sit directly after the "Corresponding author: ... Anthony
# Yoshimura" paragraph (i.e. right before the phone-number paragraph),
# instead of after the email-address paragraph where they currently live.
#
# Note: paragraph/range object references can become stale once the
# document is mutated (an insertion can "steal" a previously-held
# reference's position), so throughout this script paragraphs are always
# re-fetched from $d.Paragraphs by freshly computed index rather than by
# reusing a variable captured before a document edit.

$d = $word.ActiveDocument

function Get-ParagraphIndexByText($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $candidate = $d.Paragraphs.Item($i)
        $trimmed = $candidate.Range.Text.TrimEnd([char]13, [char]7)
        if ($trimmed -eq $text) {
            return $i
        }
    }
    return $null
}

# 1) Insert two new blank paragraphs immediately before the phone-number
#    paragraph. InsertParagraphBefore() on a range clones the paragraph /
#    run formatting of that paragraph (sz/szCs 21, spacing after=0
#    line=240 auto), matching the formatting the moved text already has.
$phoneIndex = Get-ParagraphIndexByText("(805) 886-4810")

$d.Paragraphs.Item($phoneIndex).Range.InsertParagraphBefore() | Out-Null
$d.Paragraphs.Item($phoneIndex + 1).Range.InsertParagraphBefore() | Out-Null

# 2) Fill in the two new blank paragraphs with the address text (the phone
#    paragraph - and everything after it - shifted down by two, so the new
#    blanks are exactly at $phoneIndex and $phoneIndex + 1).
$d.Paragraphs.Item($phoneIndex).Range.Text = "4875 Primrose Ln"
$d.Paragraphs.Item($phoneIndex + 1).Range.Text = "Livermore, CA 94551"

# 3) Remove the original two address paragraphs, which still exist further
#    down, immediately after the email-address paragraph.
$emailIndex = Get-ParagraphIndexByText("yoshimura4@llnl.gov")
$oldAddrIndex = $emailIndex + 1

$deleteStart = $d.Paragraphs.Item($oldAddrIndex).Range.Start
$deleteEnd = $d.Paragraphs.Item($oldAddrIndex + 1).Range.End
$d.Range($deleteStart, $deleteEnd).Delete() | Out-Null
